$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset Yousef Haidari's row back to "absent" since no departure was recorded.
$ws.Range("B5").Value = "absent"
$ws.Range("D5").Value = "15:35:13"
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
